# Apply the "minor corrections to one slide and one method" commit.
#
# 1) Re-number one slide's internal <p:sldId id="..."> (396 -> 417) for the
#    66th slide in the deck (the "Example: Parsing a Procedure Declaration"
#    slide) -- this is a purely internal bookkeeping id, exposed via the
#    Slide's SlideID property.
#
# 2) On that same slide, the Java-like pseudocode in the big code listing is
#    edited: trailing statement semicolons are dropped (the class was
#    evidently moved from a semicolon-terminated pseudo-language to one
#    without), "Token token = scanner.getToken();" becomes
#    "val token = scanner.token", and "scanner.getSymbol().isParameterDeclStarter"
#    becomes "scanner.symbol.isParameterDeclStarter".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Part 1: slide id renumbering (396 -> 417) for the slide at position 66.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(66)
$slide.SlideID = 417

# ---------------------------------------------------------------------------
# Part 2: code-listing text corrections on the "Content Placeholder 2" shape.
# ---------------------------------------------------------------------------
$sh = $slide.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

function Replace-Chars {
    param($Start, $Length, $Expected, $NewText)
    $rng = $tr.Characters($Start, $Length)
    if ($rng.Text -ne $Expected) {
        Write-Host "WARNING: at $Start,$Length expected [$Expected] but found [$($rng.Text)]"
    }
    $rng.Text = $NewText
}

# Work from the end of the text range backwards so that earlier (lower)
# character offsets stay valid while later ones are being rewritten.

# "    parseStatements();" -> "    parseStatements()"
Replace-Chars 373 22 "    parseStatements();" "    parseStatements()"

# "    parseInitialDecls();" -> "    parseInitialDecls()"
Replace-Chars 348 24 "    parseInitialDecls();" "    parseInitialDecls()"

# "    match(Symbol.leftBrace);" -> "    match(Symbol.leftBrace)"
Replace-Chars 345 2 ");" ")"

# "    match(Symbol.rightParen);" -> "    match(Symbol.rightParen)"
Replace-Chars 316 2 ");" ")"

# "        parseFormalParameters();" -> "        parseFormalParameters()"
Replace-Chars 284 3 "();" "()"

# "    if (scanner.getSymbol().isParameterDeclStarter())"
#   -> "    if (scanner.symbol.isParameterDeclStarter())"
# Two runs change inside this paragraph; edit the later one first.
Replace-Chars 229 22 "isParameterDeclStarter" "scanner.symbol.isParameterDeclStarter"
Replace-Chars 201 28 "    if (scanner.getSymbol()." "    if ("

# "    idTable.openScope(ScopeLevel.LOCAL);" -> "...LOCAL)"
Replace-Chars 197 2 ");" ")"

# "match(Symbol.leftParen);" -> "match(Symbol.leftParen)"
Replace-Chars 125 24 "match(Symbol.leftParen);" "match(Symbol.leftParen)"

# "idTable.add(procId, IdType.procedureId);" -> "...procedureId)"
Replace-Chars 102 22 ", IdType.procedureId);" ", IdType.procedureId)"

# "match(Symbol.identifier);" -> "match(Symbol.identifier)"
Replace-Chars 58 25 "match(Symbol.identifier);" "match(Symbol.identifier)"

# Paragraph 2: "Token procId = scanner.getToken();"
#   -> "val procId = scanner.token"
# Broken into separate runs (matching the authored edit): "val" / " " /
# "procId" (untouched) / " = " / "scanner.token".

# Rewrite the " = scanner.getToken();" tail run's text in one shot, then
# split "scanner.token" back off into its own run (distinct from the
# " = " run that precedes it).
Replace-Chars 35 22 " = scanner.getToken();" " = scanner.token"
$newTailRng = $tr.Characters(38, 13)
if ($newTailRng.Text -ne "scanner.token") {
    Write-Host "WARNING: tail split range got [$($newTailRng.Text)]"
}
$newTailRng.Font.Bold = $true
$newTailRng.Font.Bold = $false

# "Token " -> " " (keep as its own run, "procId" run in between is untouched)
Replace-Chars 23 6 "Token " " "

# Insert the new leading "val" run before the (now single) leading space.
$leadRng = $tr.Characters(23, 1)
if ($leadRng.Text -ne " ") {
    Write-Host "WARNING: lead range got [$($leadRng.Text)]"
}
$leadRng.Text = "val "
$valRng = $tr.Characters(23, 3)
$valRng.Font.Bold = $true
$valRng.Font.Bold = $false

# "match(Symbol.procRW);" -> "match(Symbol.procRW)"
Replace-Chars 20 2 ");" ")"

Write-Host "Final shape text:"
Write-Host $tr.Text
